# Apply results-documentation updates to the PerformanceTable worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the previously-empty Efficiency columns (Q:S) for rows 25-27.
$ws.Range("Q25").Value = 1.0400908449642901
$ws.Range("R25").Value = 2.02898199572072
$ws.Range("S25").Value = 0.874733921825536

$ws.Range("Q26").Value = 0.87382680949817904
$ws.Range("R26").Value = 1.74552186749098
$ws.Range("S26").Value = 0.65473294332200205

$ws.Range("Q27").Value = 0.75470235591996904
$ws.Range("R27").Value = 1.2281261172402
$ws.Range("S27").Value = 0.638434277152512

# Add the new Efficiency values for rows 28-29 (cells did not previously exist).
$ws.Range("Q28").Value = 0.69075647124898698
$ws.Range("R28").Value = 0.97493184796435794
$ws.Range("S28").Value = 0.55125293686424304

$ws.Range("Q29").Value = 0.67437209753178495
$ws.Range("R29").Value = 0.85422543687141195
$ws.Range("S29").Value = 0.48331781021119702

# Update the active selection on the sheet to match the author's final cursor position.
$ws.Range("U23").Select() | Out-Null
